$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (MAE)
$ws.Range("B2").Value = 0.214
$ws.Range("C2").Value = 0.246
$ws.Range("D2").Value = 0.173
$ws.Range("E2").Value = 0.314
$ws.Range("F2").Value = 0.256

# Row 3 (MSE)
$ws.Range("B3").Value = 0.18
$ws.Range("C3").Value = 0.171
$ws.Range("D3").Value = 0.08500000000000001
$ws.Range("E3").Value = 0.543
$ws.Range("F3").Value = 0.198

# Row 5 (mean Y-predicted)
$ws.Range("B5").Value = 18.285
$ws.Range("C5").Value = 15.301
$ws.Range("D5").Value = 13.131
$ws.Range("E5").Value = 31.133
$ws.Range("F5").Value = 18.158

# Row 6 (R2)
$ws.Range("B6").Value = 0.981
$ws.Range("C6").Value = 0.987
$ws.Range("D6").Value = 0.975
$ws.Range("E6").Value = 0.987
$ws.Range("F6").Value = 0.99
